$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / update the "through" date from May 14 to May 15
$ws.Name = "Through 2022-05-15"
$ws.Range("B1").Value = "May 2022 (through May 15)"

# Row 5 - Garfield Park
$ws.Range("AA5").Value = 2
$ws.Range("AF5").Value = 5

# Row 6 - Chicago Lawn
$ws.Range("AF6").Value = 1

# Row 7 - North Lawndale
$ws.Range("B7").Value = 3

# Row 11 - Roseland
$ws.Range("L11").Value = 2

# Row 14 - Lincoln Park
$ws.Range("B14").Value = 1

# Row 15 - Lake View
$ws.Range("L15").Value = 1

# Row 31 - Uptown
$ws.Range("G31").Value = 1
$ws.Range("AF31").Value = 1

# Row 32 - United Center
$ws.Range("B32").Value = 1

# Row 34 - South Deering
$ws.Range("V34").Value = 1

# Row 54 - Bridgeport
$ws.Range("B54").Value = 1

# Row 60 - East Side
$ws.Range("Q60").Value = 1

# Row 61 - East Village
$ws.Range("G61").Value = 2

# Row 72 - Lower West Side
$ws.Range("V72").Value = 1
